$wb = $excel.ActiveWorkbook

$wsProps = $wb.Worksheets.Item("Properties")

# Add the new "enabled" column to the Properties sheet
$wsProps.Range("C1").Value = "enabled"
$wsProps.Range("C2").Value = $true
$wsProps.Range("C3").Value = $true

# Update the selection on the Properties sheet
$wsProps.Range("C4").Select()

# Make "Properties" the active sheet/tab
$wsProps.Activate()
